# Adds "longest" / "shortest" example-length columns to the basic-stats
# table. Six new columns are inserted between the existing "*_avg_tokens"
# columns (G) and the "*_hapaxes" columns (old H:M, now shifted to N:S),
# which also pushes the "*_unknown" columns (old N:S) to T:Y.
#
# Inserting whole columns (rather than writing into a wider range) makes
# Excel itself perform that shift for every row - including the header
# row - so we only need to fill in the six brand-new header labels and
# the six brand-new per-language data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank columns at H:M, pushing everything from the old H
# column onward six places to the right (H->N, ... S->Y).
$ws.Range("H1:M1").EntireColumn.Insert()

# New header labels for the inserted columns.
$newHeaders = @{
    "H1" = "train_longest"
    "I1" = "dev_longest"
    "J1" = "test_longest"
    "K1" = "train_shortest"
    "L1" = "dev_shortest"
    "M1" = "test_shortest"
}
foreach ($addr in $newHeaders.Keys) {
    $ws.Range($addr).Value = $newHeaders[$addr]
}

# New per-row values: train_longest(H), dev_longest(I), test_longest(J),
# train_shortest(K), dev_shortest(L), test_shortest(M).
# Rows 9 and 10 only have "test" split data, so only J/M are populated
# there, matching the rest of those (sparse) rows.
$longestShortest = @{
    2  = @{ "H"=151; "I"=153; "J"=128; "K"=1; "L"=1;  "M"=1 }
    3  = @{ "H"=187; "I"=165; "J"=194; "K"=1; "L"=3;  "M"=2 }
    4  = @{ "H"=97;  "I"=93;  "J"=92;  "K"=2; "L"=2;  "M"=2 }
    5  = @{ "H"=159; "I"=124; "J"=106; "K"=1; "L"=1;  "M"=2 }
    6  = @{ "H"=258; "I"=195; "J"=268; "K"=3; "L"=4;  "M"=3 }
    7  = @{ "H"=166; "I"=132; "J"=153; "K"=6; "L"=11; "M"=9 }
    8  = @{ "H"=58;  "I"=46;  "J"=42;  "K"=2; "L"=3;  "M"=3 }
    9  = @{ "J"=180; "M"=7 }
    10 = @{ "J"=102; "M"=2 }
    11 = @{ "H"=256; "I"=226; "J"=238; "K"=2; "L"=3;  "M"=4 }
    12 = @{ "H"=552; "I"=463; "J"=170; "K"=1; "L"=2;  "M"=1 }
    13 = @{ "H"=118; "I"=89;  "J"=81;  "K"=4; "L"=4;  "M"=4 }
    14 = @{ "H"=82;  "I"=78;  "J"=83;  "K"=2; "L"=3;  "M"=3 }
    15 = @{ "H"=183; "I"=138; "J"=203; "K"=3; "L"=4;  "M"=3 }
    16 = @{ "H"=136; "I"=131; "J"=94;  "K"=1; "L"=1;  "M"=1 }
    17 = @{ "H"=932; "I"=263; "J"=582; "K"=1; "L"=3;  "M"=2 }
    18 = @{ "H"=283; "I"=229; "J"=222; "K"=1; "L"=3;  "M"=4 }
    19 = @{ "H"=185; "I"=126; "J"=133; "K"=3; "L"=5;  "M"=4 }
    20 = @{ "H"=213; "I"=152; "J"=157; "K"=2; "L"=3;  "M"=3 }
}

foreach ($row in $longestShortest.Keys) {
    $rowData = $longestShortest[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
